# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets
# F2: 1817 -> 1822
# F3: 8226 -> 8247
# F5: 315  -> 319

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 1822
    $ws.Range("F3").Value = 8247
    $ws.Range("F5").Value = 319
}
